$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal")

# --- 1. New data rows 13-16 (Story001 / Story002 Analysis & Conception entries) ---
$ws.Range("A13").Value = 43593
$ws.Range("B13").Value = 0.75
$ws.Range("C13").Value = "Analyse"
$ws.Range("D13").Value = "Analyse des modifications à apportées à l'API pour répondre aux critères du cahier des charges. J'ai pu remarquer que le cahier des charges intègre de nouvelles informations aux activités tel que le lieu et pays de réalisation de l'activité."

$ws.Range("A14").Value = 43593
$ws.Range("B14").Value = 1
$ws.Range("C14").Value = "Conception"
$ws.Range("D14").Value = "Je mets à jour le MCD et le MLD pour correspondre aux attentes du cahier des charges. J'ajoute ensuite au rapport de travail le nouveau MCD et le nouveau MLD"

$ws.Range("A15").Value = 43593
$ws.Range("B15").Value = 1
$ws.Range("C15").Value = "Implémentation"
$ws.Range("D15").Value = "Je mets à jour le schéma de base de données à l'aide de MySQL Workbecnh. Je créer un script permettant d'insérer toutes les localités de suisse dans la base de données. J'ai trouvé un fichier contenant les localités sur le site admin.ch"

$ws.Range("A16").Value = 43593
$ws.Range("B16").Value = 1
$ws.Range("C16").Value = "Implémentation"
$ws.Range("D16").Value = "J'enregistre le nouveau endpoint (Country) sur l'API, j'ai ensuite créer les méthodes permettant de sélectionner les pays et les localités de ceux-ci"
$ws.Range("E16").Value = "Je me suis aperçu que le format de retour de l'API n'est pas tout à fait adapté quand beaucoup d'enregistrement sont envoyés."

# --- 2. Apply the wrapped-text style (same as D4/D11 "Description" cells) to the
#        Description (D13-D16) and Retour d'expérience (E11-E45) columns, matching the
#        wrap-text formatting used throughout the rest of the table ---
$ws.Range("D4").Copy() | Out-Null
$ws.Range("D13:D16").PasteSpecial(-4122) | Out-Null
$ws.Range("E11:E45").PasteSpecial(-4122) | Out-Null

# --- 3. Extend template (empty) rows so the journal covers rows 17-45, reusing the
#        existing "date" (A) and "wrap" (D) styles from the existing template rows ---
$ws.Range("A17").Copy() | Out-Null
$ws.Range("A17:A45").PasteSpecial(-4122) | Out-Null

$ws.Range("D18").Copy() | Out-Null
$ws.Range("D25:D45").PasteSpecial(-4122) | Out-Null

# D24 switches from the "applyNumberFormat" wrap style to the plain wrap style
$ws.Range("D4").Copy() | Out-Null
$ws.Range("D24").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- 4. Row heights for the new content rows ---
$ws.Rows.Item(10).RowHeight = 30
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(14).RowHeight = 45
$ws.Rows.Item(15).RowHeight = 60
$ws.Rows.Item(16).RowHeight = 60

# --- 5. Column width tweaks ---
$ws.Columns.Item(3).ColumnWidth = 14.5703125
$ws.Columns.Item(4).ColumnWidth = 57.140625

# --- 6. Resize the "Journal" table & its AutoFilter to cover the new rows ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E45"))

# --- 7. Extend the Type (C) column data-validation list down to row 45 ---
$ws.Range("C2:C45").Validation.Delete()
$ws.Range("C2:C45").Validation.Add(3, 1, 1, "Gestion, Analyse, Conception, Implémentation, ")
$ws.Range("C2:C45").Validation.InCellDropdown = $true
$ws.Range("C2:C45").Validation.IgnoreBlank = $true
$ws.Range("C2:C45").Validation.ShowInput = $true
$ws.Range("C2:C45").Validation.ShowError = $true

# --- 8. Scroll position / selection, matching the author's final view ---
$ws.Application.Goto($ws.Range("A14"))
$ws.Range("C26:C27").Select()
